# Daily attendance processing - 2025-11-10 05:23:45
# Reorders the "Recorded By" (column G) values so that "System"/"system"
# entries are listed before the associated email address(es).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Text

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "backup@backdoor.com, System, system") {
        $cell.Value = "backup@backdoor.com, system, System"
    }
}
